$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(8,3).Value = "Gestão de Estoques"
$ws.Cells.Item(8,4).Value = "Tipos de estoque para <b>Ballou (2006)</b>"
$ws.Cells.Item(8,5).Value = "<ul>`n`t<li>Estoque de Segurança ou estoque isolado</li>`n`t<li>Estoque cíclico</li>`n`t<li>Estoque de desacoplamento</li>`n`t<li>Estoque de antecipação</li>`n`t<li>Estoques no canal de distribuição</li>`n</ul>"
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(9,3).Value = "Gestão de Estoques"
$ws.Cells.Item(9,4).Value = "tipos de estoque para <b>Brandalise (2017)</b>"
$ws.Cells.Item(9,5).Value = "<ul>`n`t<li>de matérias primas</li>`n`t<li>de produtos em processo</li>`n`t<li>de produtos acabados</li>`n`t<li>de peças de manutenção ou de reposição</li>`n</ul>"
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 0

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(10,3).Value = "Gestão de Estoques"
$ws.Cells.Item(10,4).Value = "Classificação dos Estoques para <b>Chiavenato (2005)</b>"
$ws.Cells.Item(10,5).Value = "<ul>`n`t<li>Estoques de matéria-prima</li>`n`t<li>Estoques de materiais em processamento</li>`n`t<li>Estoques de materiais semiacabados</li>`n`t<li>Estoques de materiais acabados</li>`n`t<li>Estoques de produtos acabados</li>`n</ul>"
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 0

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(11,3).Value = "Gestão de Estoques"
$ws.Cells.Item(11,4).Value = "Métodos de controle de estoques <b>(Chiavenato - 2006)</b>"
$ws.Cells.Item(11,5).Value = "<ul>`n`t<li>Sistema de duas gavetas;</li>`n`t<li>Sistema dos máximos-mínimos;</li>`n`t<li>Sistema de reposições periódicas; e</li>`n`t<li>MRP.</li>`n</ul>`n"
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 0

$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(12,3).Value = "Gestão de Estoques"
$ws.Cells.Item(12,4).Value = "Métodos de Avaliação (e controle) de Estoques"
$ws.Cells.Item(12,5).Value = "<ul>`n`t<li>LIFO</li>`n`t<li>FIFO</li>`n`t<li>Custo Médio</li>`n`t<li>Custo de Reposição.</li>`n</ul>"
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 0

$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(13,3).Value = "Gestão de Estoques"
$ws.Cells.Item(13,4).Value = "Fórmula do LEC"
$ws.Cells.Item(13,5).Value = "<b>Q = sqrt((2xDxCp)/Ce)</b>`n<ul>`n`t<li>Q = Lote</li>`n`t<li>D = Demanada</li>`n`t<li>Cp = Custo de Pedido Unitário</li>`n`t<li>Ce = Custo de manutenção de armazenamento</li>`n</ul>"
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0

$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(14,3).Value = "Gestão de Estoques"
$ws.Cells.Item(14,4).Value = "Limitações da LEC (Segundo Dias, 2010)"
$ws.Cells.Item(14,5).Value = "<ol>`n`t<li>Admite que;</li>`n`t<li>Custos de Armazenagem;</li>`n`t<li>O Custo poderá;</li>`n`t<li>Não leva em consideração;</li>`n`t<li>A fórmula se baseia;</li>`n`t<li>Algumas empresas existem; e</li>`n`t<li>Taxa de reabastecimento</li>`n</ol>"
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 0

$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(15,3).Value = "Engenharia de Métodos e Processos"
$ws.Cells.Item(15,4).Value = "gráficos ou recursos esquemáticos mais usualmente utilizados na Engenharia de Métodos e Processos"
$ws.Cells.Item(15,5).Value = "<ul>`n`t<li>Gráfico do fluxo do processo;</li>`n`t<li>Mapofluxograma;</li>`n`t<li>Carta de-para;</li>`n`t<li>Gráfico homem-máquina;</li>`n`t<li>Gráfico das duas mãos.</li>`n</ul>"
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 0

$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(16,3).Value = "Engenharia de Métodos e Processos"
$ws.Cells.Item(16,4).Value = "principais sistemas de medição do trabalho:"
$ws.Cells.Item(16,5).Value = "<ul>`n`t<li>cronometragem</li>`n`t<li>amostragem do trabalho</li>`n</ul>"
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 0

$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(17,3).Value = "Engenharia de Métodos e Processos"
$ws.Cells.Item(17,4).Value = "sistema de medição do trabalho: <b>cronometragem</b>"
$ws.Cells.Item(17,5).Value = "<ol>`n`t<li>Obter e registrar</li>`n`t<li>Dividir a operação</li>`n`t<li>Observar e registrar</li>`n`t<li>Determinar o número</li>`n`t<li>Avaliar o ritmo do operador.</li>`n`t<li>Determinar o tempo normal.</li>`n`t<li>Determinar as tolerâncias.</li>`n`t<li>Determinar o tempo-padrão da operação.</li>`n</ol>"
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 0

$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(18,3).Value = "Engenharia de Métodos e Processos"
$ws.Cells.Item(18,4).Value = "sistema de medição de trabalho: <b>amostragem do trabalho</b>"
$ws.Cells.Item(18,5).Value = "estimativa da proporção de tempo gasto em uma determinada atividade, durante um certo período, por meio de observações instantâneas, ininterruptas e espaçadas aleatoriamente."
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 0
